$d = $word.ActiveDocument

# Unicode helpers - the PS string-literal lexer here folds "smart" typographic
# punctuation (en dash, curly quotes, ...) down to plain ASCII, so build the
# handful of special characters we need from their code points instead.
$rsquo  = [char]8217   # U+2019 RIGHT SINGLE QUOTATION MARK
$endash = [char]8211   # U+2013 EN DASH

# ---------------------------------------------------------------------------
# 1. "Napkins" -> "Visualization " + "Napkin" + "s" (three underlined runs),
#    keeping the paragraph's own (pPr) underline mark intact.
# ---------------------------------------------------------------------------
$napkinsPara = $d.Paragraphs.Item(8)
if ($napkinsPara.Range.Text.TrimEnd([char]13, [char]7) -ne "Napkins") {
    throw "Unexpected paragraph 8 text: $($napkinsPara.Range.Text)"
}
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Visualization </w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>Napkin</w:t></w:r><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>s</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$napkinsPara.Range.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------------
# 2. Drop the "Nick's:" paragraph entirely.
# ---------------------------------------------------------------------------
$nick = $d.Paragraphs.Item(9)
if ($nick.Range.Text.TrimEnd([char]13, [char]7) -ne ("Nick" + $rsquo + "s:")) {
    throw "Unexpected paragraph 9 text: $($nick.Range.Text)"
}
$nick.Range.Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3. After "Vahram's:" (now paragraph 9) and before "Isabella's:", add five
#    new feedback paragraphs.
# ---------------------------------------------------------------------------
$vahram = $d.Paragraphs.Item(9)
if ($vahram.Range.Text.TrimEnd([char]13, [char]7) -ne ("Vahram" + $rsquo + "s:")) {
    throw "Unexpected paragraph 9 text: $($vahram.Range.Text)"
}

$newLines = @(
    "Scatter Plot 1 - Change Atkinson" + $rsquo + "s index to Gini Index",
    "Scatter Plot 2 " + $endash + " Everything looks fine",
    "Map " + $endash + " Too fancy for the presentation, might be difficult to absorb the information",
    "Timeline " + $endash + " Connect the two timelines to be one single timeline, also we can have 2 axis for y axis (for example, income on the right side and index on the left side)",
    "Histogram " + $endash + " Change the number of bins if necessary, do not leave it to Excel, because it might not accurately represent the picture"
)

$anchor = $vahram
foreach ($line in $newLines) {
    $anchor.Range.InsertParagraphAfter() | Out-Null
    $anchor = $d.Paragraphs.Item($anchor.Index + 1)
    $anchor.Range.Text = $line
}

# ---------------------------------------------------------------------------
# 4. After "Isabella's:" add three more paragraphs: a plain note, an
#    underlined "Presentation Napkin" heading (run-level underline only, no
#    paragraph-mark formatting), and the closing feedback paragraph.
# ---------------------------------------------------------------------------
$isabella = $d.Paragraphs.Item($anchor.Index + 1)
$isabellaPrefix = "Isabella" + $rsquo + "s:"
if ($isabella.Range.Text.Substring(0, $isabellaPrefix.Length) -ne $isabellaPrefix) {
    throw "Unexpected Isabella paragraph text: $($isabella.Range.Text)"
}

$isabella.Range.InsertParagraphAfter() | Out-Null
$pythonPara = $d.Paragraphs.Item($isabella.Index + 1)
$pythonPara.Range.Text = "Include Python and PowerBI visualizations."

$pythonPara.Range.InsertParagraphAfter() | Out-Null
$presNapkin = $d.Paragraphs.Item($pythonPara.Index + 1)
$presNapkin.Range.Text = "Presentation Napkin"

$presNapkin.Range.InsertParagraphAfter() | Out-Null
$closingPara = $d.Paragraphs.Item($presNapkin.Index + 1)
$closingPara.Range.Text = "Restating/reminding the audience of the questions is a good practice. Make sure the axes on graphs are clearly labeled. Split the questions into two slides to avoid clutter; make sure there is not too much text/clutter in general."

# Underline just the "Presentation Napkin" text (not its trailing paragraph
# mark) as the very last step, once every paragraph mark already exists -
# doing it earlier leaks the underline into the next InsertParagraphAfter().
$textOnly = $d.Range($presNapkin.Range.Start, $presNapkin.Range.End - 1)
$textOnly.Font.Underline = 1

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
